$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row values (Date serial, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
# taken from the reshuffled source rows, keyed by destination row number.
$rowData = @{
    2  = @(44763, 50,  2300, 2300, 2300, 2300)
    3  = @(44435, 130, 1300, 1300, 1300, 1300)
    4  = @(44357, 35,  1000, 1000, 1000, 1000)
    5  = @(44749, 120, 2300, 2300, 2300, 2300)
    6  = @(44748, 300, 2300, 2300, 2300, 2300)
    7  = @(44431, 100, 1300, 1300, 1300, 1300)
    8  = @(44762, 50,  2300, 2300, 2300, 2300)
    10 = @(44343, 60,  1300, 1300, 1300, 1300)
    11 = @(44438, 60,  1200, 1200, 1200, 1200)
    12 = @(44811, 60,  2500, 2500, 2500, 2500)
    13 = @(44753, 160, 2300, 2300, 2300, 2300)
    14 = @(44473, 120, 1200, 1200, 1200, 1200)
    15 = @(44424, 50,  1200, 1200, 1200, 1200)
    16 = @(44760, 80,  2300, 2300, 2300, 2300)
    17 = @(44405, 50,  1200, 1200, 1200, 1200)
    18 = @(44830, 50,  2500, 2500, 2500, 2500)
    19 = @(44432, 30,  1300, 1300, 1300, 1300)
    20 = @(44476, 80,  1200, 1200, 1200, 1200)
    21 = @(44418, 40,  1200, 1200, 1200, 1200)
    22 = @(44812, 50,  2500, 2500, 2500, 2500)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("M$r").Value = $vals[1]
    $ws.Range("N$r").Value = $vals[2]
    $ws.Range("O$r").Value = $vals[3]
    $ws.Range("P$r").Value = $vals[4]
    $ws.Range("S$r").Value = $vals[5]
}
